$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT into a target cell without altering the
# target cells existing style (used for D-column numbers that would
# otherwise be auto-coerced to the Number type by plain Value assignment).
function Set-TextValue($targetAddr, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163, $null, $false, $false)
    $scratch.Delete(-4159)
}

$ws.Range("D2").Value = "27.234.00"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "1.564.81"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue "D5" "210.81"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  -0.31%  "
Set-TextValue "D11" "0.0871"
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("D12").Value = "1.787.39"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "1.566.68"
$ws.Range("E13").Value = "  +0.54%  "
Set-TextValue "D14" "3.76"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "27.203.30"
Set-TextValue "D17" "61.94"
$ws.Range("E17").Value = "  +0.04%  "
Set-TextValue "D18" "218.02"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "0.0₃0701"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  -0.14%  "
Set-TextValue "D22" "4.14"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("E24").Value = "  +0.71%  "
Set-TextValue "D25" "151.48"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  +1.54%  "
Set-TextValue "D28" "15.03"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("D33").Value = "1.458.94"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  +1.18%  "
Set-TextValue "D40" "5.89"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("E43").Value = "  +1.38%  "
Set-TextValue "D44" "0.987"
$ws.Range("E44").Value = "  -1.62%  "
Set-TextValue "D45" "64.47"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "1.700.21"
$ws.Range("E47").Value = "  +0.36%  "
Set-TextValue "D48" "85.81"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("E50").Value = "  +1.23%  "
Set-TextValue "D51" "0.0949"
$ws.Range("E51").Value = "  -1.16%  "

$excel.CutCopyMode = $false
Write-Host "Applied cryptos update"
